{"js": "// Replace the date heading and each \"a\u00f7b=c, r\" answer cell text in place.\n// Pairs are [oldText, newText] applied in document order; every oldText is\n// a unique, exact paragraph/cell text, so a literal (non-wildcard) search\n// unambiguously matches exactly one run and insertText(..., replace) swaps\n// just that run's text while leaving its formatting (font/size) untouched.\nconst pairs = [\n  [\"2025-09-08 Monday\", \"2025-09-09 Tuesday\"],\n  [\"299\u00f76=49, 5\", \"272\u00f77=38, 6\"],\n  [\"459\u00f76=76, 3\", \"800\u00f76=133, 2\"],\n  [\"513\u00f77=73, 2\", \"193\u00f76=32, 1\"],\n  [\"657\u00f78=82, 1\", \"945\u00f75=189, 0\"],\n  [\"327\u00f75=65, 2\", \"554\u00f77=79, 1\"],\n  [\"777\u00f73=259, 0\", \"386\u00f73=128, 2\"],\n  [\"178\u00f73=59, 1\", \"639\u00f79=71, 0\"],\n  [\"762\u00f74=190, 2\", \"199\u00f72=99, 1\"],\n  [\"624\u00f75=124, 4\", \"915\u00f75=183, 0\"],\n  [\"278\u00f78=34, 6\", \"723\u00f79=80, 3\"],\n  [\"365\u00f74=91, 1\", \"541\u00f72=270, 1\"],\n  [\"332\u00f73=110, 2\", \"328\u00f78=41, 0\"],\n  [\"644\u00f77=92, 0\", \"572\u00f72=286, 0\"],\n  [\"636\u00f79=70, 6\", \"467\u00f72=233, 1\"],\n  [\"775\u00f76=129, 1\", \"878\u00f76=146, 2\"],\n  [\"597\u00f75=119, 2\", \"689\u00f76=114, 5\"],\n  [\"557\u00f72=278, 1\", \"247\u00f76=41, 1\"],\n  [\"970\u00f76=161, 4\", \"100\u00f79=11, 1\"],\n  [\"675\u00f72=337, 1\", \"575\u00f79=63, 8\"],\n  [\"894\u00f78=111, 6\", \"903\u00f78=112, 7\"],\n  [\"491\u00f77=70, 1\", \"594\u00f75=118, 4\"],\n  [\"777\u00f75=155, 2\", \"822\u00f73=274, 0\"],\n  [\"744\u00f77=106, 2\", \"794\u00f78=99, 2\"],\n  [\"891\u00f78=111, 3\", \"552\u00f78=69, 0\"],\n  [\"726\u00f74=181, 2\", \"437\u00f79=48, 5\"]\n];\n\nfor (const [oldText, newText] of pairs) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length !== 1) {\n    throw new Error(`Expected exactly 1 match for \"${oldText}\", found ${results.items.length}`);\n  }\n\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# [oldText, newText] pairs, applied in document order. Every oldText is a\n# unique, literal (non-wildcard) string, so wdReplaceAll touches exactly one run.\n$pairs = @(\n    @(\"2025-09-08 Monday\", \"2025-09-09 Tuesday\"),\n    @(\"299\u00f76=49, 5\", \"272\u00f77=38, 6\"),\n    @(\"459\u00f76=76, 3\", \"800\u00f76=133, 2\"),\n    @(\"513\u00f77=73, 2\", \"193\u00f76=32, 1\"),\n    @(\"657\u00f78=82, 1\", \"945\u00f75=189, 0\"),\n    @(\"327\u00f75=65, 2\", \"554\u00f77=79, 1\"),\n    @(\"777\u00f73=259, 0\", \"386\u00f73=128, 2\"),\n    @(\"178\u00f73=59, 1\", \"639\u00f79=71, 0\"),\n    @(\"762\u00f74=190, 2\", \"199\u00f72=99, 1\"),\n    @(\"624\u00f75=124, 4\", \"915\u00f75=183, 0\"),\n    @(\"278\u00f78=34, 6\", \"723\u00f79=80, 3\"),\n    @(\"365\u00f74=91, 1\", \"541\u00f72=270, 1\"),\n    @(\"332\u00f73=110, 2\", \"328\u00f78=41, 0\"),\n    @(\"644\u00f77=92, 0\", \"572\u00f72=286, 0\"),\n    @(\"636\u00f79=70, 6\", \"467\u00f72=233, 1\"),\n    @(\"775\u00f76=129, 1\", \"878\u00f76=146, 2\"),\n    @(\"597\u00f75=119, 2\", \"689\u00f76=114, 5\"),\n    @(\"557\u00f72=278, 1\", \"247\u00f76=41, 1\"),\n    @(\"970\u00f76=161, 4\", \"100\u00f79=11, 1\"),\n    @(\"675\u00f72=337, 1\", \"575\u00f79=63, 8\"),\n    @(\"894\u00f78=111, 6\", \"903\u00f78=112, 7\"),\n    @(\"491\u00f77=70, 1\", \"594\u00f75=118, 4\"),\n    @(\"777\u00f75=155, 2\", \"822\u00f73=274, 0\"),\n    @(\"744\u00f77=106, 2\", \"794\u00f78=99, 2\"),\n    @(\"891\u00f78=111, 3\", \"552\u00f78=69, 0\"),\n    @(\"726\u00f74=181, 2\", \"437\u00f79=48, 5\"),\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $matched = $find.Execute(\n        $oldText,   # FindText\n        $true,      # MatchCase\n        $false,     # MatchWholeWord\n        $false,     # MatchWildcards\n        $false,     # MatchSoundsLike\n        $false,     # MatchAllWordForms\n        $true,      # Forward\n        1,          # Wrap: wdFindContinue\n        $false,     # Format\n        $newText,   # ReplaceWith\n        2           # Replace: wdReplaceAll\n    )\n    if (-not $matched) {\n        throw \"Find/Replace did not match: $oldText\"\n    }\n}\n"}
